$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text-preserving format, then restore clean (no-style) text cells
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '64.119.69'
Set-TextValue 'E2' '  -3.90%  '
Set-TextValue 'D3' '3.143.67'
Set-TextValue 'E3' '  -3.51%  '
Set-TextValue 'E4' '  +0.25%  '
Set-TextValue 'D5' '607.50'
Set-TextValue 'E5' '  -0.10%  '
Set-TextValue 'D6' '146.72'
Set-TextValue 'E6' '  -7.25%  '
Set-TextValue 'E7' '  +0.16%  '
Set-TextValue 'D8' '3.136.08'
Set-TextValue 'E8' '  -3.69%  '
Set-TextValue 'D9' '0.523'
Set-TextValue 'E9' '  -4.93%  '
Set-TextValue 'D10' '0.150'
Set-TextValue 'E10' '  -7.07%  '
Set-TextValue 'D11' '5.47'
Set-TextValue 'E11' '  -7.00%  '
Set-TextValue 'D12' '0.474'
Set-TextValue 'E12' '  -6.35%  '
Set-TextValue 'D13' '0.0000249'
Set-TextValue 'E13' '  -8.51%  '
Set-TextValue 'D14' '35.41'
Set-TextValue 'E14' '  -9.80%  '
Set-TextValue 'D15' '3.663.60'
Set-TextValue 'E15' '  -3.42%  '
Set-TextValue 'D16' '64.134.78'
Set-TextValue 'E16' '  -4.01%  '
Set-TextValue 'E17' '  +0.23%  '
Set-TextValue 'D18' '3.151.03'
Set-TextValue 'E18' '  -3.24%  '
Set-TextValue 'D19' '6.88'
Set-TextValue 'E19' '  -7.53%  '
Set-TextValue 'D20' '478.56'
Set-TextValue 'E20' '  -6.14%  '
Set-TextValue 'D21' '14.71'
Set-TextValue 'E21' '  -4.76%  '
Set-TextValue 'D22' '0.710'
Set-TextValue 'E22' '  -5.74%  '
Set-TextValue 'D23' '7.74'
Set-TextValue 'E23' '  -4.90%  '
Set-TextValue 'D24' '13.58'
Set-TextValue 'E24' '  -8.73%  '
Set-TextValue 'D25' '83.25'
Set-TextValue 'E25' '  -4.04%  '
Set-TextValue 'E26' '  -0.17%  '
Set-TextValue 'D27' '2.87'
Set-TextValue 'E27' '  -5.46%  '
Set-TextValue 'D28' '8.40'
Set-TextValue 'E28' '  -7.97%  '
Set-TextValue 'D29' '2.17'
Set-TextValue 'E29' '  -10.05%  '
Set-TextValue 'D30' '6.75'
Set-TextValue 'E30' '  -1.83%  '
Set-TextValue 'D31' '0.113'
Set-TextValue 'E31' '  -20.04%  '
Set-TextValue 'B32' 'FirstDigitalUSD'
Set-TextValue 'C32' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D32' '1.00'
Set-TextValue 'E32' '  +0.11%  '
Set-TextValue 'B33' 'Stacks'
Set-TextValue 'C33' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D33' '2.72'
Set-TextValue 'E33' '  -6.81%  '
Set-TextValue 'D34' '26.09'
Set-TextValue 'E34' '  -7.26%  '
Set-TextValue 'D35' '1.09'
Set-TextValue 'E35' '  -5.31%  '
Set-TextValue 'D36' '54.14'
Set-TextValue 'E36' '  -2.85%  '
Set-TextValue 'D37' '5.93'
Set-TextValue 'E37' '  -8.59%  '
Set-TextValue 'D38' '0.0₃0727'
Set-TextValue 'E38' '  -8.91%  '
Set-TextValue 'D39' '459.63'
Set-TextValue 'E39' '  -7.24%  '
Set-TextValue 'D40' '2.92'
Set-TextValue 'E40' '  -14.34%  '
Set-TextValue 'D41' '0.0394'
Set-TextValue 'E41' '  -7.84%  '
Set-TextValue 'D42' '0.118'
Set-TextValue 'E42' '  -8.49%  '
Set-TextValue 'D43' '8.39'
Set-TextValue 'E43' '  -5.27%  '
Set-TextValue 'D44' '2.841.69'
Set-TextValue 'E44' '  -4.49%  '
Set-TextValue 'D45' '0.264'
Set-TextValue 'E45' '  -10.40%  '
Set-TextValue 'D46' '2.24'
Set-TextValue 'E46' '  -10.91%  '
Set-TextValue 'D47' '26.43'
Set-TextValue 'E47' '  -8.65%  '
Set-TextValue 'E48' '  -0.04%  '
Set-TextValue 'D49' '2.30'
Set-TextValue 'E49' '  -7.68%  '
Set-TextValue 'E50' '  -5.25%  '
Set-TextValue 'D51' '119.62'
Set-TextValue 'E51' '  -1.64%  '
